# Add 2022-Q3 data
#
# 1. Shift the existing summary rows in "总计" down by one and insert the
#    new 2022-Q3 summary row at the top of the data (row 2).
# 2. Insert a brand-new worksheet named "2022-Q3" right after "总计",
#    built by duplicating the "2021-Q4" sheet (so it inherits matching
#    styles/column layout) and then overwriting its data with the
#    2022-Q3 fund holdings.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$val) {
    # Force the cell to be treated as text even when the value looks like
    # a number (e.g. "009658" or "0.85"), then drop back to the default
    # "Normal" style so we don't leave a stray number-format style behind.
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: push rows 2-4 down to 3-5, and write the new
#    2022-Q3 row into row 2.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Give new row 5 the same style as row 4 before filling it in.
$summary.Range("A4").Copy()
$summary.Range("A5").PasteSpecial(-4122)

$summary.Range("B5").Value = $summary.Range("B4").Value()
$summary.Range("C5").Value = $summary.Range("C4").Value()
$summary.Range("D5").Value = "0.03"
$summary.Range("A5").Value = 3

$summary.Range("B4").Value = $summary.Range("B3").Value()
$summary.Range("C4").Value = $summary.Range("C3").Value()
$summary.Range("D4").Value = "0.08"
$summary.Range("A4").Value = 2

$summary.Range("B3").Value = $summary.Range("B2").Value()
$summary.Range("C3").Value = $summary.Range("C2").Value()
$summary.Range("D3").Value = "0.28"
$summary.Range("A3").Value = 1

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = "0.02"
$summary.Range("A2").Value = 0

# ---------------------------------------------------------------------
# 2) New "2022-Q3" worksheet, placed right after "总计".
#    Duplicate "2021-Q4" to inherit its column styling, rename it, then
#    extend it to 4 data rows and fill in the 2022-Q3 fund holdings.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($null, $summary)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The template only has 2 data rows (rows 2-3); extend formatting down to
# rows 4 and 5 by copying row 3's formats.
$q3.Range("A3:H3").Copy()
$q3.Range("A4:H4").PasteSpecial(-4122)
$q3.Range("A3:H3").Copy()
$q3.Range("A5:H5").PasteSpecial(-4122)

# Row 2
$q3.Range("A2").Value = 0
Set-TextValue $q3.Range("B2") "009658"
Set-TextValue $q3.Range("C2") "汇丰晋信中小盘低波动策略股票A"
Set-TextValue $q3.Range("D2") "0.85"
Set-TextValue $q3.Range("E2") "90.14"
Set-TextValue $q3.Range("F2") "1.93"
Set-TextValue $q3.Range("G2") "0.0164"
$q3.Range("H2").Value = 6

# Row 3
$q3.Range("A3").Value = 1
Set-TextValue $q3.Range("B3") "013802"
Set-TextValue $q3.Range("C3") "财通资管中证钢铁指数A"
Set-TextValue $q3.Range("D3") "0.08"
Set-TextValue $q3.Range("E3") "92.45"
Set-TextValue $q3.Range("F3") "3.25"
Set-TextValue $q3.Range("G3") "0.0026"
$q3.Range("H3").Value = 7

# Row 4
$q3.Range("A4").Value = 2
Set-TextValue $q3.Range("B4") "009775"
Set-TextValue $q3.Range("C4") "汇丰晋信中小盘低波动策略股票C"
Set-TextValue $q3.Range("D4") "0.04"
Set-TextValue $q3.Range("E4") "90.14"
Set-TextValue $q3.Range("F4") "1.93"
Set-TextValue $q3.Range("G4") "0.0008"
$q3.Range("H4").Value = 6

# Row 5
$q3.Range("A5").Value = 3
Set-TextValue $q3.Range("B5") "013803"
Set-TextValue $q3.Range("C5") "财通资管中证钢铁指数C"
Set-TextValue $q3.Range("D5") "0.01"
Set-TextValue $q3.Range("E5") "92.45"
Set-TextValue $q3.Range("F5") "3.25"
Set-TextValue $q3.Range("G5") "0.0003"
$q3.Range("H5").Value = 7

# Copying a sheet makes the copy the active tab; restore the original
# active sheet ("2020-Q4") so tabSelected stays where it was.
$wb.Worksheets.Item("2020-Q4").Activate()
